$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A24").Value = 10
$ws.Range("B24").Value = 1.5
$ws.Range("A25").Value = 20
$ws.Range("B25").Formula = "=A25*B24/A24"

$ws.Range("C21").Select()
